$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (A1:O1) to the underscored column-name convention,
# matching the "Element_Code"/"Length_m"/... style used across the project.
$ws.Range("B1").Value = "Element_Type"
$ws.Range("A1").Value = "Element_Code"
$ws.Range("C1").Value = "Material_1"
$ws.Range("D1").Value = "Material_2"
$ws.Range("E1").Value = "Material_3"
$ws.Range("F1").Value = "Length_m"
$ws.Range("G1").Value = "Width_m"
$ws.Range("H1").Value = "Height_m"
$ws.Range("I1").Value = "Thickness_1_m"
$ws.Range("J1").Value = "Thickness_2_m"
$ws.Range("K1").Value = "Thickness_3_m"
$ws.Range("L1").Value = "Surface_m2"
$ws.Range("M1").Value = "Volume_m3"
$ws.Range("N1").Value = "Slope"
$ws.Range("O1").Value = "Azimuth"

# Column A widened slightly to fit "Element_Code"; columns D/E now share a
# common width now that "Material_2"/"Material_3" are the same length.
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(4).ColumnWidth = 8.983072916666666
$ws.Columns.Item(5).ColumnWidth = 8.983072916666666

# Move the active selection back to the top of the data.
$ws.Range("A2").Select() | Out-Null
